$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Overview sheet: status columns (zh-cn / de-de) flip from "Ready for
#    handoff" to "Handed back: in sync with en-US" for both rows.
#    Same shared string is reused on the language sheets' Status column.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: Status column + newly populated "Latest Target File" /
#    "Latest Handback File" columns (handback has just been generated).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("J2").Value = "ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.f7fc272bfe7ed118da11824e0c87e1a6e27dbeea.zh-cn.xlf"
$wsZh.Range("J3").Value = "e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.bc763140408beb65b28d61b8bf1dfa8a125877da.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.md", "", "", "ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.md", "", "", "e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.md")

# ---------------------------------------------------------------------------
# 3. de-de sheet: Status column, newly populated "Latest Target File" /
#    "Latest Handback File" / "Latest Handback DateTime" columns.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("J2").Value = "ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.f7fc272bfe7ed118da11824e0c87e1a6e27dbeea.de-de.xlf"
$wsDe.Range("J3").Value = "e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.bc763140408beb65b28d61b8bf1dfa8a125877da.de-de.xlf"

$wsDe.Range("K2").Value = "2016-08-22 15:12:00"
$wsDe.Range("K3").Value = "2016-08-22 15:12:00"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.md", "", "", "ccf55d22-0fd0-4a3b-b56f-7a1a80e1fc30.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ca12a4ddc2e5c67161feb058152740a97b3776e/e2e/e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.md", "", "", "e35b6dea-0c69-4d2d-8d58-aa7954d74b1e.md")

# ---------------------------------------------------------------------------
# 4. zh-cn / de-de sheets: the "Latest Handback DateTime" cells (K2/K3) keep
#    referencing the shared placeholder string, but that string's text was
#    globally rewritten from "0001-01-01 00:00:00" to "2016-08-22 15:11:53"
#    (a side effect of the handback run timestamping every still-placeholder
#    cell). Re-assert the zh-cn values so they pick up the new text too.
# ---------------------------------------------------------------------------
$wsZh.Range("K2").Value = "2016-08-22 15:11:53"
$wsZh.Range("K3").Value = "2016-08-22 15:11:53"

# ---------------------------------------------------------------------------
# 5. Widen columns that now hold the longer hyperlink / handback file names.
# ---------------------------------------------------------------------------
$ws1.Range("E1").ColumnWidth = 29.9777047293527
$ws1.Range("F1").ColumnWidth = 29.9777047293527

$wsZh.Range("C1").ColumnWidth = 29.9777047293527
$wsZh.Range("I1").ColumnWidth = 40
$wsZh.Range("J1").ColumnWidth = 40

$wsDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDe.Range("I1").ColumnWidth = 40
$wsDe.Range("J1").ColumnWidth = 40
